$d = $word.ActiveDocument

# Locate the paragraph that contains the text to be corrected/expanded.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("O xadrez me ensinou a pensar com calma")
if (-not $found) {
    throw "Target paragraph not found"
}
$p = $rng.Paragraphs.Item(1)
$target = $p.Range

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:eastAsia="Arial" w:cs="Arial"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR"/></w:rPr>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>'
$xml += '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
$xml += '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
$xml += '<pkg:xmlData>'
$xml += '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'

$xml += '<w:p>'
$xml += '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="360"/>' + $rPr + '</w:pPr>'
$xml += '<w:r>' + $rPr + '<w:t>O xadrez me ensinou a pensar com calma, analisar possibilidades e tomar decisões de forma mais consciente. A</w:t></w:r>'
$xml += '<w:r>' + $rPr + '<w:t xml:space="preserve">credito que o </w:t></w:r>'
$xml += '<w:r>' + $rPr + '<w:t>jogo reflet</w:t></w:r>'
$xml += '<w:r>' + $rPr + '<w:t xml:space="preserve">e </w:t></w:r>'
$xml += '<w:r>' + $rPr + '<w:t>a própria vida: às vezes é preciso recuar para avançar, planejar cada passo e aceitar que nem todas as jogadas dão certo, mas sempre há uma próxima oportunidade.</w:t></w:r>'
$xml += '</w:p>'

$xml += '<w:p>'
$xml += '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="360"/>' + $rPr + '</w:pPr>'
$xml += '</w:p>'

$xml += '</w:body></w:document>'
$xml += '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
